$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.8
$ws.Range("C3").Value = 7.09
$ws.Range("C4").Value = 5.62
$ws.Range("C5").Value = 2.45
$ws.Range("C6").Value = 2.29
$ws.Range("C7").Value = 3.39
$ws.Range("C51").Value = 4.28
